$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- Update computed values (fix bug in time computing / Greedy inconsistency) ---
$ws.Range("B1").Value = 1025.16541234155
$ws.Range("C1").Value = 9.84434383314962
$ws.Range("G1").Value = 9.84184719718607
$ws.Range("Q1").Value = 9.8584060415703
$ws.Range("T1").Value = 0.0705883886709028
$ws.Range("U1").Value = 9.93816174207165
$ws.Range("AE1").Value = 9.8584060415703
$ws.Range("AH1").Value = 7.74611733640112
$ws.Range("AI1").Value = 9.93816174207165
$ws.Range("AL1").Value = 5.30026367718394
$ws.Range("AP1").Value = 5.30026367718394
$ws.Range("B2").Value = 2050.33
$ws.Range("C2").Value = 19.762205894692976
$ws.Range("G2").Value = 19.74300441899194
$ws.Range("Q2").Value = 27.440729406901
$ws.Range("R2").Value = 7.19381562608776
$ws.Range("S2").Value = 2
$ws.Range("T2").Value = 0.171717395186256
$ws.Range("U2").Value = 27.4353533954569
$ws.Range("V2").Value = 2
$ws.Range("AE2").Value = 20.8596036137675
$ws.Range("AH2").Value = 19.7865123662125
$ws.Range("AI2").Value = 20.8021578047076
$ws.Range("AL2").Value = 9.903368526732175
$ws.Range("AP2").Value = 9.903368526732175
$ws.Range("B3").Value = 3075.495
$ws.Range("C3").Value = 29.11597225162197
$ws.Range("G3").Value = 29.085727412886488
$ws.Range("Q3").Value = 51.6799717957695
$ws.Range("R3").Value = 21.5984419019716
$ws.Range("S3").Value = 6
$ws.Range("T3").Value = 0.271830303739019
$ws.Range("U3").Value = 51.4962211613742
$ws.Range("V3").Value = 6
$ws.Range("AE3").Value = 35.3579120605517
$ws.Range("AH3").Value = 53.4974850390242
$ws.Range("AI3").Value = 35.2699876051962
$ws.Range("AL3").Value = 12.407563508323898
$ws.Range("AP3").Value = 12.407563508323898
$ws.Range("B4").Value = 4100.66
$ws.Range("C4").Value = 41.924944702572205
$ws.Range("G4").Value = 41.83229423516639
$ws.Range("Q4").Value = 76.9201943544839
$ws.Range("R4").Value = 36.0090357113469
$ws.Range("S4").Value = 10
$ws.Range("T4").Value = 0.276374300148701
$ws.Range("U4").Value = 76.8618103826891
$ws.Range("V4").Value = 10
$ws.Range("AE4").Value = 52.2129161943855
$ws.Range("AH4").Value = 63.4636800286973
$ws.Range("AI4").Value = 52.2696822268459
$ws.Range("AL4").Value = 16.86054556582374
$ws.Range("AP4").Value = 16.86054556582374
$ws.Range("B5").Value = 5125.825
$ws.Range("C5").Value = 52.81492432362033
$ws.Range("G5").Value = 52.65464066142461
$ws.Range("Q5").Value = 109.466480025275
$ws.Range("R5").Value = 57.6598248538241
$ws.Range("S5").Value = 16
$ws.Range("T5").Value = 0.466225409402146
$ws.Range("U5").Value = 109.254393402142
$ws.Range("V5").Value = 16
$ws.Range("AE5").Value = 76.9767291151596
$ws.Range("AH5").Value = 66.2134789424364
$ws.Range("AI5").Value = 76.9500177987469
$ws.Range("AL5").Value = 20.319508511507696
$ws.Range("AP5").Value = 20.319508511507696
$ws.Range("B6").Value = 6150.99
$ws.Range("C6").Value = 66.44182642568548
$ws.Range("G6").Value = 65.58673389077605
$ws.Range("Q6").Value = 134.022989343576
$ws.Range("R6").Value = 72.0947926678625
$ws.Range("S6").Value = 20
$ws.Range("T6").Value = 0.472207034305553
$ws.Range("U6").Value = 133.889766564952
$ws.Range("V6").Value = 20
$ws.Range("AE6").Value = 90.6232700270524
$ws.Range("AH6").Value = 154.299997392595
$ws.Range("AI6").Value = 90.6550430330901
$ws.Range("AL6").Value = 23.661184528677854
$ws.Range("AP6").Value = 23.661184528677854
$ws.Range("B7").Value = 7176.155
$ws.Range("C7").Value = 76.7519765455253
$ws.Range("G7").Value = 75.01164664174148
$ws.Range("Q7").Value = 165.61646794222
$ws.Range("R7").Value = 93.8065307814269
$ws.Range("S7").Value = 26
$ws.Range("T7").Value = 0.617977585795488
$ws.Range("U7").Value = 165.079967327145
$ws.Range("V7").Value = 26
$ws.Range("AE7").Value = 130.425667589637
$ws.Range("AH7").Value = 92.6651980916953
$ws.Range("AI7").Value = 130.483709908652
$ws.Range("AL7").Value = 27.182327681147452
$ws.Range("AP7").Value = 27.182327681147452
$ws.Range("B8").Value = 8201.32
$ws.Range("C8").Value = 90.48839189977792
$ws.Range("G8").Value = 89.39466256925654
$ws.Range("Q8").Value = 195.130593453625
$ws.Range("R8").Value = 111.897193828006
$ws.Range("S8").Value = 31
$ws.Range("T8").Value = 0.731525347930589
$ws.Range("U8").Value = 194.742175443099
$ws.Range("V8").Value = 31
$ws.Range("AE8").Value = 135.292783343302
$ws.Range("AH8").Value = 159.589309559064
$ws.Range("AI8").Value = 135.289994215477
$ws.Range("AL8").Value = 34.8787593296795
$ws.Range("AP8").Value = 34.8787593296795
$ws.Range("B9").Value = 9226.485
$ws.Range("C9").Value = 97.76761508148334
$ws.Range("G9").Value = 95.81786858747208
$ws.Range("Q9").Value = 243.324400927429
$ws.Range("R9").Value = 137.273276343588
$ws.Range("S9").Value = 38
$ws.Range("T9").Value = 0.853296708506304
$ws.Range("U9").Value = 242.41288566777
$ws.Range("V9").Value = 38
$ws.Range("AE9").Value = 167.597943573992
$ws.Range("AH9").Value = 150.591247730372
$ws.Range("AI9").Value = 167.902984484639
$ws.Range("AL9").Value = 42.03318046067842
$ws.Range("AP9").Value = 42.03318046067842
$ws.Range("B10").Value = 10251.65
$ws.Range("C10").Value = 109.5009830679254
$ws.Range("G10").Value = 108.291607532838
$ws.Range("Q10").Value = 532.978602813924
$ws.Range("R10").Value = 162.680124704164
$ws.Range("S10").Value = 45
$ws.Range("T10").Value = 0.778438891110012
$ws.Range("U10").Value = 532.456058361574
$ws.Range("V10").Value = 45
$ws.Range("AE10").Value = 186.8062509994592
$ws.Range("AH10").Value = 117.346218393111
$ws.Range("AI10").Value = 186.2633421396792
$ws.Range("AL10").Value = 48.945727191344986
$ws.Range("AP10").Value = 48.945727191344986

# --- Update the active window view: scroll position + selection ---
try {
    $excel.ActiveWindow.ScrollColumn = 27
    $excel.ActiveWindow.ScrollRow = 1
} catch {
    # Headless runtime may not track window scroll position; ignore.
}
[void]$ws.Range("AL15").Select()
